$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
try {
  $tcs = $sm.ThemeColorScheme
  Write-Output ("sm.ThemeColorScheme null: " + ($null -eq $tcs))
} catch {
  Write-Output ("err: " + $_.Exception.Message)
}
